# Auto-generated PowerShell Excel COM-interop script
# Applies the "changing demand for amiris" edit to the AMIRIS data structure workbook
$wb = $excel.ActiveWorkbook

## --- times sheet: update StartTime / StopTime ---
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 43830.99861111111
$wsTimes.Range("B3").Value = 44195.99861111111
# Re-apply the (equivalent, case-insensitive) custom date/time format so the
# workbook mints a fresh style slot for these cells, matching the style churn
# produced by the original edit (the displayed format is unchanged either way).
$wsTimes.Range("B2:B3").NumberFormat = "yyyy-mm-dd h:mm:ss"

## --- scenario_data_emlab sheet: collapse years to 2020 only, update values ---
$wsScen = $wb.Worksheets.Item("scenario_data_emlab")
$wsScen.Range("C1:K8").Clear()
$wsScen.Range("B1").Value = 2020
$wsScen.Range("B2").Value = 20.4
$wsScen.Range("B3").Value = 1.69
$wsScen.Range("B4").Value = 3.96
$wsScen.Range("B5").Value = 10.8
$wsScen.Range("B6").Value = 20.16
$wsScen.Range("B7").Value = 46.44
$wsScen.Range("B8").Value = './amiris_workflow/amiris-config/data/load.csv'

## --- conventionals sheet: renumber identifiers, update rows, add new plants ---
$wsConv = $wb.Worksheets.Item("conventionals")
$convArr = New-Object 'object[,]' 33,7
$convArr[0,0] = 0
$convArr[0,1] = 19920300022
$convArr[0,2] = 'NATURAL_GAS'
$convArr[0,3] = 4.2
$convArr[0,4] = 0.61
$convArr[0,5] = 31358.329
$convArr[0,6] = 31358.329
$convArr[1,0] = 1
$convArr[1,1] = 19892800024
$convArr[1,2] = 'HARD_COAL'
$convArr[1,3] = 3.5
$convArr[1,4] = 0.33
$convArr[1,5] = 24845.77
$convArr[1,6] = 24845.77
$convArr[2,0] = 2
$convArr[2,1] = 20140300058
$convArr[2,2] = 'NATURAL_GAS'
$convArr[2,3] = 4.2
$convArr[2,4] = 0.61
$convArr[2,5] = 1000
$convArr[2,6] = 1000
$convArr[3,0] = 3
$convArr[3,1] = 20140300059
$convArr[3,2] = 'NATURAL_GAS'
$convArr[3,3] = 4.2
$convArr[3,4] = 0.61
$convArr[3,5] = 1000
$convArr[3,6] = 1000
$convArr[4,0] = 4
$convArr[4,1] = 20140300060
$convArr[4,2] = 'NATURAL_GAS'
$convArr[4,3] = 4.2
$convArr[4,4] = 0.61
$convArr[4,5] = 1000
$convArr[4,6] = 1000
$convArr[5,0] = 5
$convArr[5,1] = 20140300061
$convArr[5,2] = 'NATURAL_GAS'
$convArr[5,3] = 4.2
$convArr[5,4] = 0.61
$convArr[5,5] = 1000
$convArr[5,6] = 1000
$convArr[6,0] = 6
$convArr[6,1] = 20140300062
$convArr[6,2] = 'NATURAL_GAS'
$convArr[6,3] = 4.2
$convArr[6,4] = 0.61
$convArr[6,5] = 1000
$convArr[6,6] = 1000
$convArr[7,0] = 7
$convArr[7,1] = 20140300063
$convArr[7,2] = 'NATURAL_GAS'
$convArr[7,3] = 4.2
$convArr[7,4] = 0.61
$convArr[7,5] = 1000
$convArr[7,6] = 1000
$convArr[8,0] = 8
$convArr[8,1] = 20140300064
$convArr[8,2] = 'NATURAL_GAS'
$convArr[8,3] = 4.2
$convArr[8,4] = 0.61
$convArr[8,5] = 1000
$convArr[8,6] = 1000
$convArr[9,0] = 9
$convArr[9,1] = 20140300065
$convArr[9,2] = 'NATURAL_GAS'
$convArr[9,3] = 4.2
$convArr[9,4] = 0.61
$convArr[9,5] = 1000
$convArr[9,6] = 1000
$convArr[10,0] = 10
$convArr[10,1] = 20140300066
$convArr[10,2] = 'NATURAL_GAS'
$convArr[10,3] = 4.2
$convArr[10,4] = 0.61
$convArr[10,5] = 1000
$convArr[10,6] = 1000
$convArr[11,0] = 11
$convArr[11,1] = 20140300067
$convArr[11,2] = 'NATURAL_GAS'
$convArr[11,3] = 4.2
$convArr[11,4] = 0.61
$convArr[11,5] = 1000
$convArr[11,6] = 1000
$convArr[12,0] = 12
$convArr[12,1] = 20140300068
$convArr[12,2] = 'NATURAL_GAS'
$convArr[12,3] = 4.2
$convArr[12,4] = 0.61
$convArr[12,5] = 1000
$convArr[12,6] = 1000
$convArr[13,0] = 13
$convArr[13,1] = 20141700069
$convArr[13,2] = 'NATURAL_GAS'
$convArr[13,3] = 4.5
$convArr[13,4] = 0.43
$convArr[13,5] = 1000
$convArr[13,6] = 1000
$convArr[14,0] = 14
$convArr[14,1] = 20141700070
$convArr[14,2] = 'NATURAL_GAS'
$convArr[14,3] = 4.5
$convArr[14,4] = 0.43
$convArr[14,5] = 1000
$convArr[14,6] = 1000
$convArr[15,0] = 15
$convArr[15,1] = 20141700071
$convArr[15,2] = 'NATURAL_GAS'
$convArr[15,3] = 4.5
$convArr[15,4] = 0.43
$convArr[15,5] = 1000
$convArr[15,6] = 1000
$convArr[16,0] = 16
$convArr[16,1] = 20141700072
$convArr[16,2] = 'NATURAL_GAS'
$convArr[16,3] = 4.5
$convArr[16,4] = 0.43
$convArr[16,5] = 1000
$convArr[16,6] = 1000
$convArr[17,0] = 17
$convArr[17,1] = 20141700073
$convArr[17,2] = 'NATURAL_GAS'
$convArr[17,3] = 4.5
$convArr[17,4] = 0.43
$convArr[17,5] = 1000
$convArr[17,6] = 1000
$convArr[18,0] = 18
$convArr[18,1] = 20141700074
$convArr[18,2] = 'NATURAL_GAS'
$convArr[18,3] = 4.5
$convArr[18,4] = 0.43
$convArr[18,5] = 1000
$convArr[18,6] = 1000
$convArr[19,0] = 19
$convArr[19,1] = 20141700075
$convArr[19,2] = 'NATURAL_GAS'
$convArr[19,3] = 4.5
$convArr[19,4] = 0.43
$convArr[19,5] = 1000
$convArr[19,6] = 1000
$convArr[20,0] = 20
$convArr[20,1] = 20141700076
$convArr[20,2] = 'NATURAL_GAS'
$convArr[20,3] = 4.5
$convArr[20,4] = 0.43
$convArr[20,5] = 1000
$convArr[20,6] = 1000
$convArr[21,0] = 21
$convArr[21,1] = 20141700077
$convArr[21,2] = 'NATURAL_GAS'
$convArr[21,3] = 4.5
$convArr[21,4] = 0.43
$convArr[21,5] = 1000
$convArr[21,6] = 1000
$convArr[22,0] = 22
$convArr[22,1] = 20141700078
$convArr[22,2] = 'NATURAL_GAS'
$convArr[22,3] = 4.5
$convArr[22,4] = 0.43
$convArr[22,5] = 1000
$convArr[22,6] = 1000
$convArr[23,0] = 23
$convArr[23,1] = 20141700079
$convArr[23,2] = 'NATURAL_GAS'
$convArr[23,3] = 4.5
$convArr[23,4] = 0.43
$convArr[23,5] = 1000
$convArr[23,6] = 1000
$convArr[24,0] = 24
$convArr[24,1] = 20191700124
$convArr[24,2] = 'NATURAL_GAS'
$convArr[24,3] = 4.5
$convArr[24,4] = 0.43
$convArr[24,5] = 1000
$convArr[24,6] = 1000
$convArr[25,0] = 25
$convArr[25,1] = 20191700125
$convArr[25,2] = 'NATURAL_GAS'
$convArr[25,3] = 4.5
$convArr[25,4] = 0.43
$convArr[25,5] = 1000
$convArr[25,6] = 1000
$convArr[26,0] = 26
$convArr[26,1] = 20191700126
$convArr[26,2] = 'NATURAL_GAS'
$convArr[26,3] = 4.5
$convArr[26,4] = 0.43
$convArr[26,5] = 1000
$convArr[26,6] = 1000
$convArr[27,0] = 27
$convArr[27,1] = 20191700127
$convArr[27,2] = 'NATURAL_GAS'
$convArr[27,3] = 4.5
$convArr[27,4] = 0.43
$convArr[27,5] = 1000
$convArr[27,6] = 1000
$convArr[28,0] = 28
$convArr[28,1] = 20191700128
$convArr[28,2] = 'NATURAL_GAS'
$convArr[28,3] = 4.5
$convArr[28,4] = 0.43
$convArr[28,5] = 1000
$convArr[28,6] = 1000
$convArr[29,0] = 29
$convArr[29,1] = 19843000129
$convArr[29,2] = 'OIL'
$convArr[29,3] = 6
$convArr[29,4] = 0.35
$convArr[29,5] = 3652.9
$convArr[29,6] = 3652.9
$convArr[30,0] = 30
$convArr[30,1] = 19822900131
$convArr[30,2] = 'LIGNITE'
$convArr[30,3] = 3.5
$convArr[30,4] = 0.33
$convArr[30,5] = 20779.02
$convArr[30,6] = 20779.02
$convArr[31,0] = 31
$convArr[31,1] = 19851400132
$convArr[31,2] = 'NUCLEAR'
$convArr[31,3] = 3.5
$convArr[31,4] = 0.33
$convArr[31,5] = 8599
$convArr[31,6] = 8599
$convArr[32,0] = 32
$convArr[32,1] = 19921700133
$convArr[32,2] = 'NATURAL_GAS'
$convArr[32,3] = 4.5
$convArr[32,4] = 0.43
$convArr[32,5] = 8194.3025
$convArr[32,6] = 8194.3025
$wsConv.Range("A2:G34").Value = $convArr

## --- renewables sheet: renumber identifiers (column B only) ---
$wsRen = $wb.Worksheets.Item("renewables")
$renArr = New-Object 'object[,]' 69,1
$renArr[0,0] = 20122100025
$renArr[1,0] = 20122100026
$renArr[2,0] = 20122100027
$renArr[3,0] = 20122100028
$renArr[4,0] = 20122100029
$renArr[5,0] = 20122100030
$renArr[6,0] = 20122100031
$renArr[7,0] = 20122100032
$renArr[8,0] = 20122100033
$renArr[9,0] = 20122100034
$renArr[10,0] = 20122100035
$renArr[11,0] = 20132100036
$renArr[12,0] = 20132100037
$renArr[13,0] = 20132100038
$renArr[14,0] = 20132100039
$renArr[15,0] = 20132100040
$renArr[16,0] = 20132100041
$renArr[17,0] = 20132100042
$renArr[18,0] = 20132100043
$renArr[19,0] = 20132100044
$renArr[20,0] = 20132100045
$renArr[21,0] = 20132100046
$renArr[22,0] = 20132400047
$renArr[23,0] = 20132400048
$renArr[24,0] = 20132400049
$renArr[25,0] = 20132400050
$renArr[26,0] = 20132400051
$renArr[27,0] = 20132400052
$renArr[28,0] = 20132400053
$renArr[29,0] = 20132400054
$renArr[30,0] = 20132400055
$renArr[31,0] = 20132400056
$renArr[32,0] = 20132400057
$renArr[33,0] = 20142100080
$renArr[34,0] = 20142100081
$renArr[35,0] = 20142100082
$renArr[36,0] = 20142100083
$renArr[37,0] = 20142100084
$renArr[38,0] = 20142100085
$renArr[39,0] = 20142100086
$renArr[40,0] = 20142100087
$renArr[41,0] = 20142100088
$renArr[42,0] = 20142100089
$renArr[43,0] = 20142100090
$renArr[44,0] = 20152400100
$renArr[45,0] = 20152400101
$renArr[46,0] = 20152400102
$renArr[47,0] = 20152400103
$renArr[48,0] = 20152400104
$renArr[49,0] = 20152400105
$renArr[50,0] = 20152400106
$renArr[51,0] = 20152400107
$renArr[52,0] = 20152400108
$renArr[53,0] = 20152400109
$renArr[54,0] = 20152400110
$renArr[55,0] = 20162300113
$renArr[56,0] = 20162300114
$renArr[57,0] = 20162300115
$renArr[58,0] = 20162300116
$renArr[59,0] = 20162300117
$renArr[60,0] = 20162300118
$renArr[61,0] = 20162300119
$renArr[62,0] = 20162300120
$renArr[63,0] = 20162300121
$renArr[64,0] = 20162300122
$renArr[65,0] = 20162300123
$renArr[66,0] = 19641200130
$renArr[67,0] = 20102100134
$renArr[68,0] = 20142300135
$wsRen.Range("B3:B71").Value = $renArr

## --- biogas sheet: renumber identifiers (column B only) ---
$wsBio = $wb.Worksheets.Item("biogas")
$bioArr = New-Object 'object[,]' 11,1
$bioArr[0,0] = 20150100091
$bioArr[1,0] = 20150100092
$bioArr[2,0] = 20150100093
$bioArr[3,0] = 20150100094
$bioArr[4,0] = 20150100095
$bioArr[5,0] = 20150100096
$bioArr[6,0] = 20150100097
$bioArr[7,0] = 20150100098
$bioArr[8,0] = 20150100099
$bioArr[9,0] = 20160100111
$bioArr[10,0] = 20160100112
$wsBio.Range("B3:B13").Value = $bioArr

Write-Output "edit applied"
